$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.478.15"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "2.777.34"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'352.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "'108.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'39.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").Value = "'20.10"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.10%  "
$ws.Range("D13").Value = "'0.0835"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").Value = "'7.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "3.215.06"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").Value = "2.764.14"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "51.500.60"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "'7.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("E20").Value = "  -0.50%  "
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "0.0₃0962"
$ws.Range("E22").Value = "  -1.69%  "
$ws.Range("D23").Value = "'69.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").Value = "'266.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.89%  "
$ws.Range("D25").Value = "'2.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'26.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.29%  "
$ws.Range("D28").Value = "'0.162"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.08%  "
$ws.Range("D29").Value = "'10.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "'2.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'36.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.04%  "
$ws.Range("D32").Value = "'6.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.70%  "
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").Value = "'0.0828"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'18.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("E39").Value = "  -2.66%  "
$ws.Range("E40").Value = "  -1.75%  "
$ws.Range("E41").Value = "  -0.88%  "
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'120.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("D44").Value = "'21.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.30%  "
$ws.Range("E45").Value = "  -2.83%  "
$ws.Range("D46").Value = "2.112.40"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("E47").Value = "  +0.71%  "
$ws.Range("E48").Value = "  +6.66%  "
$ws.Range("E49").Value = "  -4.95%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("E51").Value = "  +8.38%  "
